$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "first-stage" / "second-stage" labels between row 2 and row 15,
# and relabel the row-15 header as "binary variables".
$ws.Range("A2").Value2 = "second-stage"
$ws.Range("A15").Value2 = "first-stage"
$ws.Range("B15").Value2 = "binary variables"

# Rows 3-8 (now the second-stage flow variables) get a "flow: " prefix.
$ws.Range("A3").Value2 = "flow: x[Arc(""s1"", ""p1"")]"
$ws.Range("A4").Value2 = "flow: x[Arc(""s2"", ""p1"")]"
$ws.Range("A5").Value2 = "flow: x[Arc(""s3"", ""t1"")]"
$ws.Range("A6").Value2 = "flow: x[Arc(""s3"", ""t2"")]"
$ws.Range("A7").Value2 = "flow: x[Arc(""p1"", ""t1"")]"
$ws.Range("A8").Value2 = "flow: x[Arc(""p1"", ""t2"")]"

# Rows 16-22 (now the first-stage binary variables) get a "decision on: " prefix.
$ws.Range("A16").Value2 = "decision on: x[Arc(""s1"", ""p1"")]"
$ws.Range("A17").Value2 = "decision on: x[Arc(""s2"", ""p1"")]"
$ws.Range("A18").Value2 = "decision on: x[Arc(""s3"", ""t1"")]"
$ws.Range("A19").Value2 = "decision on: x[Arc(""s3"", ""t2"")]"
$ws.Range("A20").Value2 = "decision on: x[Arc(""p1"", ""t1"")]"
$ws.Range("A21").Value2 = "decision on: x[Arc(""p1"", ""t2"")]"
$ws.Range("A22").Value2 = "decision on: p1"
